$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Activate()

# Write the brand-new shared strings first, in the exact order they were
# first introduced by the original author, so the regenerated shared
# string table lines up with the target workbook.
$ws.Range("D41").Value = "test slug"
$ws.Range("B43").Value = "tip test"
$ws.Range("G43").Value = "tc"
$ws.Range("H41").Value = " test c"
$ws.Range("G35").Value = " content"
$ws.Range("G37").Value = "ttc "
$ws.Range("H37").Value = " abc"
$ws.Range("I43").Value = " test"

# Remaining cells that reuse already-existing shared strings (order among
# these does not affect shared string allocation).
$ws.Range("C33").Value = "test tip"
$ws.Range("C37").Value = "testtip"
$ws.Range("E37").Value = "desc"
$ws.Range("C39").Value = "test tip"

$ws.Range("B41").Value = "test tip"
$ws.Range("C41").Value = 12
$ws.Range("E41").Value = "test"
$ws.Range("F41").Value = "desc"
$ws.Range("G41").Value = "ttc"

$ws.Range("C43").Value = 12
$ws.Range("D43").Value = "test slug"
$ws.Range("E43").Value = "test tip"
$ws.Range("F43").Value = "desc"
$ws.Range("H43").Value = "ttc"
$ws.Range("J43").Value = " content"

# Column G width (nearest value the host's column-width rounding can
# represent; target stored width is 15.85546875)
$ws.Columns.Item(7).ColumnWidth = 15

# Sheet view: top left cell and selection
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("I43").Select()

$wb.Save()
